$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 693.6667
$ws.Range("I5").Value = 39.5
$ws.Range("J5").Value = 2002
$ws.Range("K5").Value = 39.5
$ws.Range("L5").Value = 2002
$ws.Range("M5").Value = 75.5
$ws.Range("N5").Value = -2232
$ws.Range("H53").Value = 79.5
$ws.Range("I53").Value = 70
$ws.Range("J53").Value = 89
$ws.Range("K53").Value = 70
$ws.Range("L53").Value = 89
$ws.Range("M53").Value = 567
$ws.Range("N53").Value = -1363
$ws.Range("H55").Value = 909.75
$ws.Range("I55").Value = 1055.6
$ws.Range("J55").Value = 805.5714
$ws.Range("K55").Value = 1055.6
$ws.Range("L55").Value = 805.5714
$ws.Range("M55").Value = -841.5999999999999
$ws.Range("N55").Value = -1233.5714
$ws.Range("H103").Value = 400
$ws.Range("I103").Value = 400
$ws.Range("J103").Value = 400
$ws.Range("K103").Value = 1200
$ws.Range("L103").Value = 1200
$ws.Range("M103").Value = -614
$ws.Range("N103").Value = -2372

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 172.5
$ws.Range("I4").Value = 172.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 172.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -56.5
$ws.Range("H32").Value = 202682.03
$ws.Range("I32").Value = 212861.56
$ws.Range("J32").Value = 141604.88
$ws.Range("K32").Value = 212861.56
$ws.Range("L32").Value = 141604.88
$ws.Range("M32").Value = -212574.56
$ws.Range("N32").Value = -142178.88
$ws.Range("H49").Value = 50026.668
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 50026.668
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 50026.668
$ws.Range("N49").Value = -50546.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 29700
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 29700
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 29700
$ws.Range("N13").Value = -30036
$ws.Range("H26").Value = 10333.333
$ws.Range("I26").Value = 10333.333
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 10333.333
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -10041.333
$ws.Range("N26").ClearContents()
$ws.Range("H86").Value = 501879
$ws.Range("I86").Value = 2027.8182
$ws.Range("J86").Value = 2334666.8
$ws.Range("K86").Value = 2027.8182
$ws.Range("L86").Value = 2334666.8
$ws.Range("M86").Value = -904.8181999999999
$ws.Range("N86").Value = -2336912.8
$ws.Range("H89").Value = 501879
$ws.Range("I89").Value = 2027.8182
$ws.Range("J89").Value = 2334666.8
$ws.Range("K89").Value = 10139.091
$ws.Range("L89").Value = 11673334
$ws.Range("M89").Value = -4523.091
$ws.Range("N89").Value = -11684566
$ws.Range("H96").Value = 16500
$ws.Range("I96").Value = 14750
$ws.Range("J96").Value = 20000
$ws.Range("K96").Value = 14750
$ws.Range("L96").Value = 20000
$ws.Range("M96").Value = -12004
$ws.Range("H105").Value = 1180.25
$ws.Range("I105").Value = 1060.5
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 1060.5
$ws.Range("L105").Value = 1300
$ws.Range("M105").Value = 686.5
$ws.Range("N105").Value = -4794
$ws.Range("H134").Value = 26371816
$ws.Range("I134").Value = 31252154
$ws.Range("J134").Value = 343338
$ws.Range("K134").Value = 93756462
$ws.Range("L134").Value = 1030014
$ws.Range("M134").Value = -93753927

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23867.043
$ws.Range("I31").Value = 30280.646
$ws.Range("J31").Value = 5695.1665
$ws.Range("K31").Value = 30280.646
$ws.Range("L31").Value = 5695.1665
$ws.Range("M31").Value = -29985.646
$ws.Range("N31").Value = -6285.1665
$ws.Range("H34").Value = 23867.043
$ws.Range("I34").Value = 30280.646
$ws.Range("J34").Value = 5695.1665
$ws.Range("K34").Value = 30280.646
$ws.Range("L34").Value = 5695.1665
$ws.Range("M34").Value = -30078.646
$ws.Range("N34").Value = -6099.1665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2233.3333
$ws.Range("I80").Value = 1300
$ws.Range("J80").Value = 2700
$ws.Range("K80").Value = 3900
$ws.Range("L80").Value = 8100
$ws.Range("M80").Value = -2964
$ws.Range("N80").Value = -9972
$ws.Range("H83").Value = 2233.3333
$ws.Range("I83").Value = 1300
$ws.Range("J83").Value = 2700
$ws.Range("K83").Value = 11700
$ws.Range("L83").Value = 24300
$ws.Range("M83").Value = -7020
$ws.Range("N83").Value = -33660
$ws.Range("H99").Value = 6149.3335
$ws.Range("I99").Value = 400
$ws.Range("J99").Value = 6868
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 20604
$ws.Range("M99").Value = 1046
$ws.Range("N99").Value = -25096

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 746.4761999999999
$ws.Range("I22").Value = 822.9091
$ws.Range("J22").Value = 662.4
$ws.Range("K22").Value = 822.9091
$ws.Range("L22").Value = 662.4
$ws.Range("M22").Value = -527.9091
$ws.Range("N22").Value = -1252.4
$ws.Range("H27").Value = 746.4761999999999
$ws.Range("I27").Value = 822.9091
$ws.Range("J27").Value = 662.4
$ws.Range("K27").Value = 822.9091
$ws.Range("L27").Value = 662.4
$ws.Range("M27").Value = -715.9091
$ws.Range("N27").Value = -876.4
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()
$ws.Range("H74").Value = 49939
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 49939
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 49939
$ws.Range("N74").Value = -51935
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 49939
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 49939
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 149817
$ws.Range("N77").Value = -159801
$ws.Range("M77").ClearContents()
$ws.Range("H93").Value = 1126.4117
$ws.Range("I93").Value = 775.75
$ws.Range("J93").Value = 1968
$ws.Range("K93").Value = 775.75
$ws.Range("L93").Value = 1968
$ws.Range("M93").Value = 472.25
$ws.Range("N93").Value = -4464
$ws.Range("H140").Value = 52857.145
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 52857.145
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 52857.145
$ws.Range("N140").Value = -63217.145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 23886.666
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 23886.666
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 23886.666
$ws.Range("N58").Value = -24502.666
$ws.Range("H122").Value = 4182.0356
$ws.Range("I122").Value = 3343.5625
$ws.Range("J122").Value = 5300
$ws.Range("K122").Value = 10030.6875
$ws.Range("L122").Value = 15900
$ws.Range("M122").Value = -7580.6875
$ws.Range("N122").Value = -20800

